# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.396.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.005"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3759"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3411"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07655"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.963"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.907"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.573.11"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.37"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06757"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.55%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.218"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5299"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.443"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.396.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.740"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.20%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.24"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.47%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.060"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.83%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.93"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.754.37"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.019"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.52%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.165"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.22%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.81%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08551"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02550"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.35%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2316"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.80%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06524"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.327"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.73%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.430"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6449"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.80%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.004"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.37%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6020"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.788"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.42%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.293"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.088"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.38"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.57%  "

Write-Output "Applied cryptos update"